$d = $word.ActiveDocument

# Locate the end of the document content (after the last paragraph,
# which currently ends with "...áttértem a pythonra.")
$endRange = $d.Range($d.Content.End, $d.Content.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" +
  '<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>5. hét</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>Működésre bírtam pythonból az UR5 kart a szimulátorral, meg is mozgattam. Ehhez készítettem egy MoveIt configot az assistanttal.</w:t></w:r></w:p>' +
  '</w:body></w:wordDocument>'

$endRange.InsertXML($xml)
